$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Extend the used range of the Metadata sheet from A1:B15 to A1:B16 ---
# Copy the formatting of the last existing data row (15) down onto the new row (16)
# so the new row picks up the same cell style without Excel fabricating an extra
# unused style entry (which happens with Rows.Insert()).
$ws.Range("A15:B15").Copy()
$ws.Range("A16:B16").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- Update simple property values (rows 1-9 are unchanged) ---
$ws.Range("B3").Value = "0.1.7"                                 # Version
$ws.Range("B6").Value = "draft"                                 # Status
$ws.Range("B8").Value = "2024-08-27T12:23:18-05:00"             # Date

# --- Rewrite rows 10-16 with the new Contact/Jurisdiction/Description/Purpose/
#     Copyright/Immutable block ---
$ws.Range("A10").Value = "Contact"
$ws.Range("B10").Value = "The Medical College of Wisconsin, Inc. and the National Marrow Donor Program (http://www.cibmtr.org)"

$ws.Range("A11").Value = "Contact"
$ws.Range("B11").Value = "Bob Milius (bmilius@nmdp.org)"

$ws.Range("A12").Value = "Jurisdiction"
$ws.Range("B12").Value = ""
# Force B12 to be stored as an (empty) shared string rather than a truly blank
# cell, matching the source workbook's representation: reuse the existing
# empty-string cell found on the "Include from LOINC" sheet (A77) as the
# source of a value-only paste so the cell keeps its text type.
$ws3 = $wb.Worksheets.Item("Include from LOINC")
$ws3.Range("A77").Copy()
$ws.Range("B12").PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$ws.Range("A13").Value = "Description"
$ws.Range("B13").Value = "Manual Differential panel - Blood (24318-8)"

$ws.Range("A14").Value = "Purpose"
$ws.Range("B14").Value = ""

$ws.Range("A15").Value = "Copyright"
$ws.Range("B15").Value = ""

$ws.Range("A16").Value = "Immutable"
$ws.Range("B16").Value = "BooleanType[null]"

Write-Host "Metadata sheet updated"
